$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting from the row above (row 15, "Bom" cell style) onto the
# new row so the added row matches the look of the rest of the table.
$ws.Range("A15:E15").Copy() | Out-Null
$ws.Range("A16:E16").PasteSpecial(-4122) | Out-Null

# New passive skill entry: Force of Habit
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "ForceOfHabit"
$ws.Range("C16").Value = "PassiveSkill"
$ws.Range("D16").Value = 0
$ws.Range("E16").Value = 0

# Remove the old empty placeholder cell that used to live at F16
$ws.Range("F16").Clear() | Out-Null

# Move the active selection down to F17, mirroring where the cursor ended up
# after the new row was inserted
$ws.Range("F17").Select() | Out-Null
